# Big 12 and updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: school names (rows 45-54), written in row order ---
$ws.Range("A45").Value = "Baylor"
$ws.Range("A46").Value = "Iowa State"
$ws.Range("A47").Value = "Kansas"
$ws.Range("A48").Value = "Kansas State"
$ws.Range("A49").Value = "Oklahoma"
$ws.Range("A50").Value = "Oklahoma State"
$ws.Range("A51").Value = "TCU"
$ws.Range("A52").Value = "Texas"
$ws.Range("A53").Value = "Texas Tech"
$ws.Range("A54").Value = "West Virginia"

# --- Column B: conference (all "Big 12") ---
$ws.Range("B45").Value = "Big 12"
$ws.Range("B46").Value = "Big 12"
$ws.Range("B47").Value = "Big 12"
$ws.Range("B48").Value = "Big 12"
$ws.Range("B49").Value = "Big 12"
$ws.Range("B50").Value = "Big 12"
$ws.Range("B51").Value = "Big 12"
$ws.Range("B52").Value = "Big 12"
$ws.Range("B53").Value = "Big 12"
$ws.Range("B54").Value = "Big 12"

# --- Column C: filenames. These were originally authored in two batches
# (single-word school scripts first, then underscored multi-word scripts),
# which is reflected in the shared-string insertion order.
$ws.Range("C45").Value = "baylor.py"
$ws.Range("C46").Value = "iowa state.py"
$ws.Range("C47").Value = "kansas.py"
$ws.Range("C49").Value = "oklahoma.py"
$ws.Range("C51").Value = "tcu.py"
$ws.Range("C52").Value = "texas.py"
$ws.Range("C48").Value = "kansas_state.py"
$ws.Range("C50").Value = "oklahoma_state.py"
$ws.Range("C53").Value = "texas_tech.py"
$ws.Range("C54").Value = "west_virginia.py"

# --- Columns D/E/F: grid/table/ul flags ---
$ws.Range("D45").Value = 1
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 0

$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 1
$ws.Range("F46").Value = 0

$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 1
$ws.Range("F47").Value = 0

$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 1
$ws.Range("F48").Value = 0

$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 1

$ws.Range("D50").Value = 1
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0

$ws.Range("D51").Value = 1
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0

$ws.Range("D52").Value = 0
$ws.Range("E52").Value = 1
$ws.Range("F52").Value = 0

$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 1
$ws.Range("F53").Value = 0

$ws.Range("D54").Value = 0
$ws.Range("E54").Value = 1
$ws.Range("F54").Value = 0

# --- Column G: scrape_date. Write raw Excel serial-date numbers (instead of
# a DateTime value) so the engine does not auto-create an extra ad-hoc date
# number format/style; the real m/d/yyyy format is copied in from the row
# above further down.
$ws.Range("G45").Value = 43211
$ws.Range("G46").Value = 43211
$ws.Range("G47").Value = 43211
$ws.Range("G48").Value = 43211
$ws.Range("G49").Value = 43211
$ws.Range("G50").Value = 43211
$ws.Range("G51").Value = 43211
$ws.Range("G52").Value = 43212
$ws.Range("G53").Value = 43212
$ws.Range("G54").Value = 43212

# Copy the date number format from the row above down onto the new G cells
$ws.Range("G44").Copy()
$ws.Range("G45:G54").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column H: football_url (hyperlinks), written in row order ---
$ws.Hyperlinks.Add($ws.Range("H45"), "http://www.baylorbears.com/sports/m-footbl/mtt/bay-m-footbl-mtt.html")
$ws.Range("H45").Value = "http://www.baylorbears.com/sports/m-footbl/mtt/bay-m-footbl-mtt.html "

$ws.Hyperlinks.Add($ws.Range("H46"), "http://cyclones.com/roster.aspx?path=football")
$ws.Range("H46").Value = "http://cyclones.com/roster.aspx?path=football "

$ws.Hyperlinks.Add($ws.Range("H47"), "https://kuathletics.com/roster.aspx?path=football")
$ws.Range("H47").Value = "https://kuathletics.com/roster.aspx?path=football "

$ws.Hyperlinks.Add($ws.Range("H48"), "http://www.kstatesports.com/roster.aspx?path=football")
$ws.Range("H48").Value = "http://www.kstatesports.com/roster.aspx?path=football "

$ws.Hyperlinks.Add($ws.Range("H49"), "http://www.soonersports.com/SportSelect.dbml?DB_OEM_ID=31000&SPID=127245&SPSID=750326")
$ws.Range("H49").Value = "http://www.soonersports.com/SportSelect.dbml?DB_OEM_ID=31000&SPID=127245&SPSID=750326 "

$ws.Hyperlinks.Add($ws.Range("H50"), "http://okstate.com/roster.aspx?path=football")
$ws.Range("H50").Value = "http://okstate.com/roster.aspx?path=football "

$ws.Hyperlinks.Add($ws.Range("H51"), "http://www.gofrogs.com/sports/m-footbl/mtt/tcu-m-footbl-mtt.html")
$ws.Range("H51").Value = "http://www.gofrogs.com/sports/m-footbl/mtt/tcu-m-footbl-mtt.html "

$ws.Hyperlinks.Add($ws.Range("H52"), "http://texassports.com/roster.aspx?path=football")
$ws.Range("H52").Value = "http://texassports.com/roster.aspx?path=football "

$ws.Hyperlinks.Add($ws.Range("H53"), "https://texastech.com/roster.aspx?path=football")
$ws.Range("H53").Value = "https://texastech.com/roster.aspx?path=football "

$ws.Hyperlinks.Add($ws.Range("H54"), "https://wvusports.com/roster.aspx?path=football")
$ws.Range("H54").Value = "https://wvusports.com/roster.aspx?path=football "

# Copy the hyperlink cell style from the row above down onto the new H cells
$ws.Range("H44").Copy()
$ws.Range("H45:H54").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- View: move the frozen pane / selection down to the new bottom of the table ---
$ws.Range("F55").Select()
